$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2000500
$ws.Range("I74").Value = 2000500
$ws.Range("K74").Value = 2000500
$ws.Range("M74").Value = -1999564

$ws.Range("H77").Value = 2000500
$ws.Range("I77").Value = 2000500
$ws.Range("K77").Value = 10002500
$ws.Range("M77").Value = -9997820

$ws.Range("H132").Value = 887.95654
$ws.Range("I132").Value = 887.95654
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2663.86962
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -133.8696199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3399
$ws.Range("I45").Value = 1997
$ws.Range("K45").Value = 1997
$ws.Range("M45").Value = -1620

$ws.Range("H61").Value = 983.3333
$ws.Range("I61").Value = 983.3333
$ws.Range("K61").Value = 983.3333
$ws.Range("M61").Value = -771.3333

$ws.Range("H74").Value = 575.3
$ws.Range("I74").Value = 417.33334
$ws.Range("K74").Value = 417.33334
$ws.Range("M74").Value = 456.66666

$ws.Range("H77").Value = 575.3
$ws.Range("I77").Value = 417.33334
$ws.Range("K77").Value = 2086.6667
$ws.Range("M77").Value = 2281.3333

$ws.Range("H122").Value = 14421.77
$ws.Range("I122").Value = 7187.5454
$ws.Range("K122").Value = 21562.6362
$ws.Range("M122").Value = -19112.6362

$ws.Range("H132").Value = 1532.4286
$ws.Range("I132").Value = 1532.4286
$ws.Range("K132").Value = 4597.2858
$ws.Range("M132").Value = -2067.2858

$ws.Range("H136").Value = 983.3333
$ws.Range("I136").Value = 983.3333
$ws.Range("K136").Value = 2949.9999
$ws.Range("M136").Value = -399.9998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 280.4
$ws.Range("I80").Value = 239.4
$ws.Range("J80").Value = 321.4
$ws.Range("K80").Value = 239.4
$ws.Range("L80").Value = 321.4
$ws.Range("M80").Value = 758.6
$ws.Range("N80").Value = -2317.4

$ws.Range("H83").Value = 280.4
$ws.Range("I83").Value = 239.4
$ws.Range("J83").Value = 321.4
$ws.Range("K83").Value = 1197
$ws.Range("L83").Value = 1607
$ws.Range("M83").Value = 3795
$ws.Range("N83").Value = -11591

$ws.Range("H86").Value = 2597.2307
$ws.Range("I86").Value = 2688.6667
$ws.Range("K86").Value = 2688.6667
$ws.Range("M86").Value = -1565.6667

$ws.Range("H89").Value = 2597.2307
$ws.Range("I89").Value = 2688.6667
$ws.Range("K89").Value = 13443.3335
$ws.Range("M89").Value = -7827.333500000001

$ws.Range("H94").Value = 1090.6666
$ws.Range("I94").Value = 631.96155
$ws.Range("K94").Value = 631.96155
$ws.Range("M94").Value = -180.96155

$ws.Range("H105").Value = 3971534.8
$ws.Range("I105").Value = 9262248
$ws.Range("K105").Value = 9262248
$ws.Range("M105").Value = -9260501

$ws.Range("H107").Value = 1342.2084
$ws.Range("I107").Value = 1168.9333
$ws.Range("K107").Value = 1168.9333
$ws.Range("M107").Value = 751.0667000000001

$ws.Range("H134").Value = 2220.6
$ws.Range("I134").Value = 1525.875
$ws.Range("K134").Value = 4577.625
$ws.Range("M134").Value = -2042.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 29991.5
$ws.Range("J3").Value = 29991.5
$ws.Range("L3").Value = 29991.5
$ws.Range("N3").Value = -30217.5

$ws.Range("H12").Value = 6435.6665
$ws.Range("I12").Value = 3155
$ws.Range("J12").Value = 12997
$ws.Range("K12").Value = 3155
$ws.Range("L12").Value = 12997
$ws.Range("M12").Value = -2985
$ws.Range("N12").Value = -13337

$ws.Range("H15").Value = 9214.5
$ws.Range("I15").Value = 14420
$ws.Range("J15").Value = 4009
$ws.Range("K15").Value = 14420
$ws.Range("L15").Value = 4009
$ws.Range("M15").Value = -14250
$ws.Range("N15").Value = -4349

$ws.Range("H105").Value = 3023.077
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 3108.3333
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 3108.3333
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6602.3333

$ws.Range("H132").Value = 1007.2222
$ws.Range("I132").Value = 1007.2222
$ws.Range("K132").Value = 3021.6666
$ws.Range("M132").Value = -491.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 88.166664
$ws.Range("I6").Value = 65.8
$ws.Range("K6").Value = 197.4
$ws.Range("M6").Value = -84.39999999999998

$ws.Range("H9").Value = 7000175
$ws.Range("I9").Value = 350
$ws.Range("K9").Value = 1050
$ws.Range("M9").Value = -826

$ws.Range("H120").Value = 9605
$ws.Range("I120").Value = 1525
$ws.Range("J120").Value = 11625
$ws.Range("K120").Value = 4575
$ws.Range("L120").Value = 34875
$ws.Range("M120").Value = 263
$ws.Range("N120").Value = -44551

$ws.Range("H138").Value = 3451.8
$ws.Range("J138").Value = 3756.3333
$ws.Range("L138").Value = 11268.9999
$ws.Range("N138").Value = -21548.9999

$ws.Range("H139").Value = 1657.3334
$ws.Range("I139").Value = 1657.3334
$ws.Range("K139").Value = 4972.0002
$ws.Range("M139").Value = 167.9997999999996

$ws.Range("H140").Value = 4997.5
$ws.Range("J140").Value = 4997.5
$ws.Range("L140").Value = 14992.5
$ws.Range("N140").Value = -25352.5

$ws.Range("H141").Value = 6146.857
$ws.Range("I141").Value = 6146.857
$ws.Range("K141").Value = 18440.571
$ws.Range("M141").Value = -13260.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2120.85
$ws.Range("I102").Value = 1165.2667
$ws.Range("K102").Value = 1165.2667
$ws.Range("M102").Value = 456.7333000000001

$ws.Range("H122").Value = 37715.25
$ws.Range("I122").Value = 1709.7084
$ws.Range("K122").Value = 5129.1252
$ws.Range("M122").Value = -2679.1252

$ws.Range("H132").Value = 1685.7142
$ws.Range("I132").Value = 1438.5
$ws.Range("K132").Value = 4315.5
$ws.Range("M132").Value = -1785.5

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = -4888
$ws.Range("N3").Value = 0

$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 5000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5000
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = -4830
$ws.Range("N15").Value = 0

$ws.Range("H16").Value = 1107.0667
$ws.Range("I16").Value = 1031.2307
$ws.Range("K16").Value = 1031.2307
$ws.Range("M16").Value = -861.2307000000001

$ws.Range("H21").Value = 3810.5
$ws.Range("I21").Value = 690
$ws.Range("J21").Value = 5544.1113
$ws.Range("K21").Value = 690
$ws.Range("L21").Value = 5544.1113
$ws.Range("M21").Value = -516
$ws.Range("N21").Value = -5892.1113

$ws.Range("H23").Value = 12335
$ws.Range("I23").Value = 12335
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 12335
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -12105

$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 8000
$ws.Range("N34").Value = -8344

$ws.Range("H47").Value = 21000
$ws.Range("J47").Value = 21000
$ws.Range("L47").Value = 21000
$ws.Range("N47").Value = -21980

$ws.Range("H52").Value = 21000
$ws.Range("J52").Value = 21000
$ws.Range("L52").Value = 21000
$ws.Range("N52").Value = -21466

$ws.Range("H61").Value = 9261185
$ws.Range("I61").Value = 12347269
$ws.Range("K61").Value = 12347269
$ws.Range("M61").Value = -12347067

$ws.Range("H93").Value = 970.94116
$ws.Range("I93").Value = 957.4286
$ws.Range("K93").Value = 957.4286
$ws.Range("M93").Value = 290.5714

$ws.Range("H113").Value = 9261185
$ws.Range("I113").Value = 12347269
$ws.Range("K113").Value = 12347269
$ws.Range("M113").Value = -12345099

$ws.Range("H134").Value = 54750
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 54750
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").Value = 54750
$ws.Range("N134").Value = -64890

$ws.Range("H136").Value = 7581.5
$ws.Range("I136").Value = 6998.75
$ws.Range("J136").Value = 8164.25
$ws.Range("K136").Value = 20996.25
$ws.Range("L136").Value = 24492.75
$ws.Range("M136").Value = -18446.25
$ws.Range("N136").Value = -29592.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 18255.666
$ws.Range("I11").Value = 14999
$ws.Range("J11").Value = 19884
$ws.Range("K11").Value = 14999
$ws.Range("L11").Value = 19884
$ws.Range("M11").Value = -14857
$ws.Range("N11").Value = -20168

$ws.Range("H12").Value = 6499.5
$ws.Range("I12").Value = 5999
$ws.Range("K12").Value = 5999
$ws.Range("M12").Value = -5857

$ws.Range("H20").Value = 46666.668
$ws.Range("J20").Value = 46666.668
$ws.Range("L20").Value = 46666.668
$ws.Range("N20").Value = -47146.668
